$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.403.14"
$ws.Range("E2").Value = "  -2.67%  "
$ws.Range("D3").Value = "1.775.59"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("E4").Value = "  -0.57%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.53%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "306.26"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("E7").Value = "  +1.28%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3612"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.61%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07151"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.22%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.8386"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +1.20%  "
$ws.Range("D12").Value = "1.770.39"
$ws.Range("E12").Value = "  -3.32%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "6.452"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.76%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.256"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.33%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.06903"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.86%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.81%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "79.02"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -1.16%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000008725"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("E19").Value = "  -0.66%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "14.91"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.07%  "
$ws.Range("D21").Value = "26.414.61"
$ws.Range("E21").Value = "  -3.10%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.099"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.16%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "10.92"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("D24").Value = "1.990.44"
$ws.Range("E24").Value = "  -5.46%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "151.62"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.98%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "1.794"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -7.86%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "18.04"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.18%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "5.068"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.64%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "114.21"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +1.51%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.841"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +11.14%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.08826"
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.7282"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +1.27%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.122"
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.319"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.51%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.66%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.733"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -4.69%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.092"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.57%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.05113"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("E39").Value = "  -0.23%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.4927"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("E41").Value = "  -0.49%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "2.600"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.26%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "6.332"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +2.69%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "8.086"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "104.72"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "10.23"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.55%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.625"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +2.52%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.06174"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.18%  "
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.4445"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.20%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.713"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +3.24%  "
